$p = $ppt.ActivePresentation

# --- 1) Table on slide 6 ("SOURCES OF FINANCE") gets a new built-in table style ---
$s = $p.Slides.Item(6)
$tbl = $s.Shapes.Item(2).Table
$tbl.ApplyStyle("{EEC21EDC-561C-4E7A-98C3-30595B527D5A}")

# --- 2) Re-colour the deck's theme (was "Integral" green palette -> standard
#        "Office Theme" palette). Order matches a:clrScheme child order:
#        dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink ---
$theme = $p.SlideMaster.Theme
$scheme = $theme.ThemeColorScheme
$scheme.Colors(1).RGB  = RGB(0x00, 0x00, 0x00)   # dk1
$scheme.Colors(2).RGB  = RGB(0xFF, 0xFF, 0xFF)   # lt1
$scheme.Colors(3).RGB  = RGB(0x44, 0x54, 0x6A)   # dk2
$scheme.Colors(4).RGB  = RGB(0xE7, 0xE6, 0xE6)   # lt2
$scheme.Colors(5).RGB  = RGB(0x5B, 0x9B, 0xD5)   # accent1
$scheme.Colors(6).RGB  = RGB(0xED, 0x7D, 0x31)   # accent2
$scheme.Colors(7).RGB  = RGB(0xA5, 0xA5, 0xA5)   # accent3
$scheme.Colors(8).RGB  = RGB(0xFF, 0xC0, 0x00)   # accent4
$scheme.Colors(9).RGB  = RGB(0x44, 0x72, 0xC4)   # accent5
$scheme.Colors(10).RGB = RGB(0x70, 0xAD, 0x47)   # accent6
$scheme.Colors(11).RGB = RGB(0x05, 0x63, 0xC1)   # hlink
$scheme.Colors(12).RGB = RGB(0x95, 0x4F, 0x72)   # folHlink
